# Insert a new weekly price record as row 101 in the "Ajo" (garlic) price
# sheet. Inserting the row shifts the existing rows 101-193 down to
# rows 102-194 (and grows the used range to A1:R194), which matches the
# rest of the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 101, pushing everything below
# down by one (old row 101 becomes row 102, ..., old row 193 becomes 194).
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new record.
$ws.Range("A101").Value = 8
$ws.Range("B101").Value = "Terminal La Palmera de La Serena"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 44554
$ws.Range("E101").Value = 4
$ws.Range("F101").Value = 100112003
$ws.Range("G101").Value = "Ajo"
$ws.Range("H101").Value = "Chino"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 560
$ws.Range("K101").Value = 19000
$ws.Range("L101").Value = 19500
$ws.Range("M101").Value = 19250
$ws.Range("N101").Value = "`$/caja 10 kilos"
$ws.Range("O101").Value = "China"
$ws.Range("P101").Value = 1925
$ws.Range("Q101").Value = 10
$ws.Range("R101").Value = "Hortaliza"
